# Apply the "add carbon credit calcs" edit to enterprise-flows-all-areas.xlsx
#
# Summary of the change:
#  - Rename sheet "soil_emissions" -> "soil_n2o"
#  - Rename sheet "soil_c" -> "carbon_credits"
#  - On the carbon_credits sheet:
#      * update the source-note row (A2) to a comet-planner URL
#      * fix a typo'd flow_type value ("output??" -> "output")
#      * rename the flow_cat value ("soil c" -> "carbon_credit")
#      * add three new data rows (8, 9, 10) covering tulare_county/
#        central_valley_organic carbon-credit co2/n2o outputs
#      * widen column E a touch for the new content
#  - Make carbon_credits the active/selected sheet (it was field_ops before)

$wb = $excel.ActiveWorkbook

# --- 1. Rename the two sheets ------------------------------------------------
$wsN2o = $wb.Worksheets.Item("soil_emissions")
$wsN2o.Name = "soil_n2o"

$ws = $wb.Worksheets.Item("soil_c")
$ws.Name = "carbon_credits"

# --- 2. Update existing rows on carbon_credits ------------------------------
$ws.Range("A2").Value = "http://comet-planner-cdfahsp.com/"

$ws.Range("B7").Value = "output"
$ws.Range("C7").Value = "carbon_credit"

# --- 3. Add the new data rows ------------------------------------------------
$ws.Range("A8").Value = "tulare_county"
$ws.Range("B8").Value = "output"
$ws.Range("C8").Value = "carbon_credit"
$ws.Range("D8").Value = "n2o"

$ws.Range("A9").Value = "central_valley_organic"
$ws.Range("B9").Value = "output"
$ws.Range("C9").Value = "carbon_credit"
$ws.Range("D9").Value = "co2"

$ws.Range("A10").Value = "central_valley_organic"
$ws.Range("B10").Value = "output"
$ws.Range("C10").Value = "carbon_credit"
$ws.Range("D10").Value = "n2o"

# match the existing red "output" styling (column B) used on row 7
$ws.Range("B8:B10").Font.Color = 255

# --- 4. Widen column E slightly (room for the new notes/values) ------------
$ws.Columns.Item(5).ColumnWidth = 16.15

# --- 5. Make carbon_credits the active sheet & set its selection -----------
$ws.Activate()
$ws.Range("E10").Select()
